$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header "上繳日" in G1 is renamed to "繳息迄日"
$ws.Range("G1").Value = "繳息迄日"

# Reflect the selection change recorded in the saved view state
$ws.Range("G1").Select()
